# Append run: update "取得日時" timestamp for existing rows, insert two new
# listings (new row 5 and new row 7 in the final layout), and append two more
# rows at the bottom (new row 12, and the shifted AWS row 13).
#
# Net effect vs. the starting sheet (9 data rows, rows 2-10):
#   - every existing listing's "取得日時" becomes 2025-12-08 12:38:30
#   - 2 brand-new listings are inserted in priority order
#   - final sheet has 12 data rows (rows 2-13)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$timestamp = "2025-12-08 12:38:30"

# Final data set, in the row order / shape described by the diff.
# Each entry: row, title, category, price, deadline, url, score, skills(optional)
$data = @(
    @{Row=2;  B="【完全在宅】ChatGPT・AI活用講師募集|IT/業務支援経験者歓迎!"; C="システム開発"; D="10,000 円 ~ 20,000 円 / 固定"; E="期限情報なし"; F="https://www.lancers.jp/work/detail/5449394"; G=600; H="🔥AI,GPT"},
    @{Row=3;  B="DreamWeaver – 夢日記 + 睡眠導入 + AI分析のアプリ開発"; C="システム開発"; D="1,000 ~ 5,000 円 / 固定"; E="期限情報なし"; F="https://www.lancers.jp/work/detail/5449048"; G=370; H="🔥AI,Ai ◆開発 ◇アプリ"},
    @{Row=4;  B="大手製造業向け センサー画像解析・高画質化のR&Dを支援するAIエンジニア募集(画像生成/超解像)"; C="システム開発"; D="300,000 円 ~ 500,000 円 / 固定"; E="期限情報なし"; F="https://www.lancers.jp/work/detail/5427956"; G=310; H="🔥AI,Ai"},
    @{Row=5;  B="【TypeScript/Clasp必須】LINE WORKS連携ファイル自動保存システムのGAS開発"; C="システム開発"; D="20,000 円 ~ 50,000 円 / 固定"; E="期限情報なし"; F="https://www.lancers.jp/work/detail/5449466"; G=193; H="🔥TypeScript ◆開発"},
    @{Row=6;  B="【フリーランス募集】訪問看護向けスケジュール管理アプリ開発"; C="システム開発"; D="1,000,000 円 ~ 3,000,000 円 / 固定"; E="期限情報なし"; F="https://www.lancers.jp/work/detail/5449280"; G=135; H="◆開発 ◇アプリ"},
    @{Row=7;  B="[週2常駐] Laravel + Vue.js 基幹業務システム開発"; C="システム開発"; D="1,000,000 円 ~ 3,000,000 円 / 固定"; E="期限情報なし"; F="https://www.lancers.jp/work/detail/5449536"; G=125; H="◆開発,システム開発"},
    @{Row=8;  B="【急募】紙の伝票をWEBシステムへ自動データ入力開発"; C="システム開発"; D="300,000 円 ~ 500,000 円 / 固定"; E="期限情報なし"; F="https://www.lancers.jp/work/detail/5449142"; G=90; H="◆開発"},
    @{Row=9;  B="自動出品システムの開発"; C="システム開発"; D="100,000 円 ~ 200,000 円 / 固定"; E="期限情報なし"; F="https://www.lancers.jp/work/detail/5449232"; G=83; H="◆開発"},
    @{Row=10; B="初回 WebアプリのiOSアプリ化+IAPサブスク(2週無料)+申請"; C="システム開発"; D="500,000 円 ~ 1,000,000 円 / 固定"; E="期限情報なし"; F="https://www.lancers.jp/work/detail/5449067"; G=45; H="◇アプリ"},
    @{Row=11; B="【急募】Shopifyでの3Dカスタムシミュレーター導入設定依頼"; C="システム開発"; D="200,000 円 ~ 300,000 円 / 固定"; E="期限情報なし"; F="https://www.lancers.jp/work/detail/5449335"; G=18; H=$null},
    @{Row=12; B="【急募】LINEのLステップ構築をサポートしてくれる方"; C="システム開発"; D="20,000 円 ~ 50,000 円 / 固定"; E="期限情報なし"; F="https://www.lancers.jp/work/detail/5449657"; G=13; H=$null},
    @{Row=13; B="初回 【継続案件】AWS上でのLAMP環境構築および保守・運用サポートパートナー募集"; C="システム開発"; D="20,000 円 ~ 50,000 円 / 固定"; E="期限情報なし"; F="https://www.lancers.jp/work/detail/5449313"; G=13; H=$null}
)

foreach ($item in $data) {
    $r = $item.Row

    $ws.Cells.Item($r, 1).Value = $timestamp
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 7).Value = $item.G

    if ($item.H -ne $null) {
        $ws.Cells.Item($r, 8).Value = $item.H
    }
}

# Column F (URL) hyperlinks. Deleting a single cell's Hyperlinks collection
# clears the *whole sheet's* hyperlink collection in this environment, so
# wipe them all in one shot and then rebuild every row's link (including the
# 3 rows whose target didn't move) to keep F2:F13 fully consistent.
if ($ws.Hyperlinks.Count -gt 0) {
    $ws.Hyperlinks.Delete()
}

$links = @(
    @{Row=2;  Url="https://www.lancers.jp/work/detail/5449394"},
    @{Row=3;  Url="https://www.lancers.jp/work/detail/5449048"},
    @{Row=4;  Url="https://www.lancers.jp/work/detail/5427956"},
    @{Row=5;  Url="https://www.lancers.jp/work/detail/5449466"},
    @{Row=6;  Url="https://www.lancers.jp/work/detail/5449280"},
    @{Row=7;  Url="https://www.lancers.jp/work/detail/5449536"},
    @{Row=8;  Url="https://www.lancers.jp/work/detail/5449142"},
    @{Row=9;  Url="https://www.lancers.jp/work/detail/5449232"},
    @{Row=10; Url="https://www.lancers.jp/work/detail/5449067"},
    @{Row=11; Url="https://www.lancers.jp/work/detail/5449335"},
    @{Row=12; Url="https://www.lancers.jp/work/detail/5449657"},
    @{Row=13; Url="https://www.lancers.jp/work/detail/5449313"}
)

foreach ($link in $links) {
    $cell = $ws.Cells.Item($link.Row, 6)
    # Pass TextToDisplay explicitly: Hyperlinks.Add leaves a non-blank cell's
    # existing text alone otherwise, which would desync the visible URL text
    # from the link target for every row whose listing shifted position.
    $ws.Hyperlinks.Add($cell, $link.Url, [Type]::Missing, [Type]::Missing, $link.Url)
    $cell.Style = "Hyperlink"
}

# Column B widened by one character unit (51 -> 52 raw OOXML width).
$ws.Columns.Item(2).ColumnWidth = 51.166666666666664
